$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced into numbers by Excel (e.g. "1.000" -> 1, "42.49" -> 42.49 float).
$textCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16",
    "D17", "D18", "D19", "D20", "D21", "D22", "D25", "D27", "D28", "D29",
    "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40",
    "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.108.32"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.916.11"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "320.34"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.5068"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.4081"
$ws.Range("E8").Value = "  +4.44%  "
$ws.Range("D9").Value = "0.08357"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").Value = "42.49"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").Value = "23.89"
$ws.Range("E12").Value = "  +5.74%  "
$ws.Range("D13").Value = "6.405"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").Value = "1.906.09"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "7.238"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "92.52"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "0.00001097"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").Value = "0.06507"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "18.52"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "5.945"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "30.112.06"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "2.191"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "2.125.14"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").Value = "21.86"
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("D28").Value = "162.92"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "2.285"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("D30").Value = "128.88"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").Value = "1.148"
$ws.Range("E31").Value = "  +10.59%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "5.963"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "3.795"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "0.02459"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").Value = "5.365"
$ws.Range("E36").Value = "  +4.19%  "
$ws.Range("D37").Value = "0.06409"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "0.2156"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "0.6556"
$ws.Range("E39").Value = "  +4.85%  "
$ws.Range("D40").Value = "1.197"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "8.643"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").Value = "11.42"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "1.215"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "13.44"
$ws.Range("E44").Value = "  +5.20%  "
$ws.Range("D45").Value = "0.6093"
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").Value = "2.193"
$ws.Range("E46").Value = "  +10.80%  "
$ws.Range("D47").Value = "3.625"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "1.211"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "122.15"
$ws.Range("D50").Value = "79.08"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("D51").Value = "1.139"
$ws.Range("E51").Value = "  -1.22%  "
